$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete row 28 entirely (the "FERRITE BEAD 976 OHM AXIAL 1LN" / "or" alternative row),
# shifting all rows below it up by one.
$ws.Rows("28").Delete()

$ws.Range("H41").Select()
